$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matching")

$ws.Range("A1").Value = "Key terms"
$ws.Range("B1").Value = "Correct order of definitions"
$ws.Range("C1").Value = "Definitions"

$ws.Range("B2").Value = "C"
$ws.Range("B3").Value = "B"
$ws.Range("B4").Value = "D"
$ws.Range("B5").Value = "A"

$ws.Range("B6").Select()
